$p = $ppt.ActivePresentation

# Unhide slides 11-20 (remove show="0")
for ($i = 11; $i -le 20; $i++) {
    $s = $p.Slides.Item($i)
    $s.SlideShowTransition.Hidden = 0
}

# Update "Result is N rows / M columns" row counts on slides 15, 16, 17, 18
# Slide 15: TextBox 2, paragraph 3 -> "Result is 5,509 rows / 14 columns" -> "3,048"
$s15 = $p.Slides.Item(15)
$tb15 = $s15.Shapes.Item(2)
$para15 = $tb15.TextFrame.TextRange.Paragraphs(3)
$para15.Characters(11, 6).Text = "3,048 "

# Slide 16: TextBox 2, paragraph 4 -> "Result is 5,509 rows / 2 columns" -> "3,048"
$s16 = $p.Slides.Item(16)
$tb16 = $s16.Shapes.Item(2)
$para16 = $tb16.TextFrame.TextRange.Paragraphs(4)
$para16.Characters(11, 6).Text = "3,048 "

# Slide 17: TextBox 2, paragraph 3 -> "Result is 67 rows / 2 columns" -> "25"
$s17 = $p.Slides.Item(17)
$tb17 = $s17.Shapes.Item(2)
$para17 = $tb17.TextFrame.TextRange.Paragraphs(3)
$para17.Characters(11, 3).Text = "25 "

# Slide 18: TextBox 2, paragraph 5 -> "Result is 67 rows / 2 columns" -> "25"
$s18 = $p.Slides.Item(18)
$tb18 = $s18.Shapes.Item(2)
$para18 = $tb18.TextFrame.TextRange.Paragraphs(5)
$para18.Characters(11, 3).Text = "25 "
